$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Segment ID and code on row 2
$ws.Range("A2").Value = 15651527
$ws.Range("B2").Value = -2147370268

# Update Segment Name and Segment Description text
$ws.Range("C2").Value = "Eyeota - US 33Across - B2B - Healthcare Executives"
$ws.Range("D2").Value = "Leadership + Healthcare Executives"
